# Commit: "update mgmt course analysis english keywords, add TUM MDT"
# For this workbook the only content change is a new program row:
#   Program_choosing!A10 = "TUM_Mgmt_DigitalTech"
#   Program_choosing!B10 = "Yes"
# plus the accompanying "Yes/No" list validation on the new cell (mirrors
# the existing validation that already covers B1:B9).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A10").Value = "TUM_Mgmt_DigitalTech"
$ws.Range("B10").Value = "Yes"

$ws.Range("B10").Validation.Add(3, 1, 1, """Yes,No""")
$ws.Range("B10").Validation.IgnoreBlank = $true
$ws.Range("B10").Validation.ShowInput = $false
$ws.Range("B10").Validation.ShowError = $true
